# Apply updated probability values to Sheet1 (Northwestern_B team matrix)
# Changes reflect team matrices recalculated from games pulled March 7
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2113821138211382
$ws.Range("C2").Value = 0.5271002710027101
$ws.Range("J2").Value = 0.01490514905149051
$ws.Range("P2").Value = 0.1626016260162602
$ws.Range("S2").Value = 0.08401084010840108
$ws.Range("B3").Value = 0.009975062344139651
$ws.Range("C3").Value = 0.02493765586034913
$ws.Range("J3").Value = 0.02743142144638404
$ws.Range("P3").Value = 0.7356608478802993
$ws.Range("S3").Value = 0.2019950124688279
$ws.Range("J4").Value = 0.01041666666666667
$ws.Range("P4").Value = 0.6145833333333334
$ws.Range("S4").Value = 0.375
$ws.Range("B6").Value = 0.06888361045130641
$ws.Range("D6").Value = 0.007125890736342043
$ws.Range("F6").Value = 0.04513064133016627
$ws.Range("J6").Value = 0.2992874109263658
$ws.Range("O6").Value = 0.03087885985748218
$ws.Range("Q6").Value = 0.1353919239904988
$ws.Range("R6").Value = 0.06175771971496437
$ws.Range("S6").Value = 0.3515439429928741
$ws.Range("B7").Value = 0.1349206349206349
$ws.Range("D7").Value = 0.01587301587301587
$ws.Range("E7").Value = 0.002645502645502645
$ws.Range("F7").Value = 0.04232804232804233
$ws.Range("J7").Value = 0.1137566137566138
$ws.Range("O7").Value = 0.01851851851851852
$ws.Range("Q7").Value = 0.1613756613756614
$ws.Range("R7").Value = 0.08994708994708994
$ws.Range("S7").Value = 0.4206349206349206
$ws.Range("B8").Value = 0.1146838156484459
$ws.Range("D8").Value = 0.02143622722400857
$ws.Range("E8").Value = 0.001071811361200429
$ws.Range("F8").Value = 0.04180064308681672
$ws.Range("J8").Value = 0.1339764201500536
$ws.Range("O8").Value = 0.02572347266881029
$ws.Range("Q8").Value = 0.1661307609860664
$ws.Range("R8").Value = 0.1028938906752412
$ws.Range("S8").Value = 0.3922829581993569
$ws.Range("B9").Value = 0.141439205955335
$ws.Range("D9").Value = 0.01240694789081886
$ws.Range("F9").Value = 0.07940446650124069
$ws.Range("J9").Value = 0.1215880893300248
$ws.Range("O9").Value = 0.03473945409429281
$ws.Range("Q9").Value = 0.1464019851116625
$ws.Range("R9").Value = 0.07196029776674938
$ws.Range("S9").Value = 0.3920595533498759
$ws.Range("B10").Value = 0.1271946937182989
$ws.Range("D10").Value = 0.02614124073351541
$ws.Range("E10").Value = 0.0007803355442840422
$ws.Range("F10").Value = 0.0647678501755755
$ws.Range("J10").Value = 0.1170503316426063
$ws.Range("O10").Value = 0.02262973078423722
$ws.Range("Q10").Value = 0.2130316035895435
$ws.Range("R10").Value = 0.07647288333983612
$ws.Range("S10").Value = 0.351931330472103
$ws.Range("G11").Value = 0.1325878594249201
$ws.Range("J11").Value = 0.1150159744408946
$ws.Range("K11").Value = 0.207667731629393
$ws.Range("L11").Value = 0.5239616613418531
$ws.Range("S11").Value = 0.0207667731629393
$ws.Range("G12").Value = 0.7391304347826086
$ws.Range("J12").Value = 0.2115942028985507
$ws.Range("K12").Value = 0.005797101449275362
$ws.Range("L12").Value = 0.01739130434782609
$ws.Range("S12").Value = 0.02608695652173913
$ws.Range("G13").Value = 0.6296296296296297
$ws.Range("J13").Value = 0.2962962962962963
$ws.Range("S13").Value = 0.07407407407407407
$ws.Range("F15").Value = 0.02355460385438972
$ws.Range("H15").Value = 0.1413276231263383
$ws.Range("I15").Value = 0.07494646680942184
$ws.Range("J15").Value = 0.291220556745182
$ws.Range("K15").Value = 0.05781584582441113
$ws.Range("M15").Value = 0.008565310492505354
$ws.Range("O15").Value = 0.07066381156316917
$ws.Range("S15").Value = 0.3319057815845824
$ws.Range("F16").Value = 0.01956521739130435
$ws.Range("H16").Value = 0.1652173913043478
$ws.Range("I16").Value = 0.08695652173913043
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.1108695652173913
$ws.Range("M16").Value = 0.01304347826086956
$ws.Range("N16").Value = 0.004347826086956522
$ws.Range("O16").Value = 0.06739130434782609
$ws.Range("S16").Value = 0.1326086956521739
$ws.Range("F17").Value = 0.02531645569620253
$ws.Range("H17").Value = 0.1657077100115075
$ws.Range("I17").Value = 0.09896432681242807
$ws.Range("J17").Value = 0.4177215189873418
$ws.Range("K17").Value = 0.1001150747986191
$ws.Range("M17").Value = 0.01726121979286536
$ws.Range("O17").Value = 0.06904487917146145
$ws.Range("S17").Value = 0.1058688147295742
$ws.Range("F18").Value = 0.01061007957559682
$ws.Range("H18").Value = 0.1750663129973475
$ws.Range("I18").Value = 0.07957559681697612
$ws.Range("J18").Value = 0.3660477453580902
$ws.Range("K18").Value = 0.1273209549071618
$ws.Range("M18").Value = 0.02122015915119363
$ws.Range("O18").Value = 0.08488063660477453
$ws.Range("S18").Value = 0.1352785145888594
$ws.Range("F19").Value = 0.01494298073141958
$ws.Range("H19").Value = 0.2308297286669288
$ws.Range("I19").Value = 0.08415257569799449
$ws.Range("J19").Value = 0.3716083366103028
$ws.Range("K19").Value = 0.1053873377900118
$ws.Range("M19").Value = 0.02084152575697994
$ws.Range("N19").Value = 0.0007864726700747149
$ws.Range("O19").Value = 0.05859221392056626
$ws.Range("S19").Value = 0.1128588281557216
